$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) "Topics: ..." line - capitalise Numbers / Casting / String in the
#    topic list (and split the run the same way the authored edit did).
# -------------------------------------------------------------------
$topicsIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.StartsWith("Topics: Syntax, Comment, Variables, Data type, numbers, casting, string")) {
        $topicsIdx = $i
    }
}

if ($topicsIdx -gt 0) {
    $p = $d.Paragraphs.Item($topicsIdx)
    $r = $p.Range
    $xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="005768E4" w:rsidRPr="005768E4" w:rsidRDefault="005768E4" w:rsidP="005768E4"><w:pPr><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="005768E4"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Topics</w:t></w:r><w:r w:rsidRPr="005768E4"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Syntax, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Comment, Variables, Data type, Numbers, Casting, S</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>tring      and Booleans</w:t></w:r></w:p>
'@
    [void]$r.InsertXML($xml)
}

# -------------------------------------------------------------------
# 2) Collapse the empty paragraph that sits between "...Boolean
#    answer:" and the following "Example:" paragraph.
# -------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($k = 2; $k -le $count - 1; $k++) {
    $prevText = $d.Paragraphs.Item($k - 1).Range.Text
    $curText = $d.Paragraphs.Item($k).Range.Text
    $nextText = $d.Paragraphs.Item($k + 1).Range.Text
    if ($prevText.StartsWith("When you compare two values") -and $curText -eq [string][char]13 -and $nextText.StartsWith("Example:")) {
        [void]$d.Paragraphs.Item($k).Range.Delete()
        break
    }
}

# -------------------------------------------------------------------
# 3) Replace the trailing empty "NormalWeb" paragraph with a new
#    "Resources" section, and strip the first-line indent from the
#    final (now last) paragraph of the document.
# -------------------------------------------------------------------
$count = $d.Paragraphs.Count
$resourcesParaIdx = $count - 1

# The target paragraph is the empty, non-indented "NormalWeb" paragraph
# that sits right before the document's very last (empty) paragraph.
# Guard against drift by scanning backwards for that shape if the
# simple offset doesn't land on an empty paragraph.
$candidate = $d.Paragraphs.Item($resourcesParaIdx)
if ($candidate.Range.Text -ne [string][char]13) {
    for ($k = $count - 1; $k -ge 2; $k--) {
        $cand = $d.Paragraphs.Item($k)
        if ($cand.Range.Text -eq [string][char]13 -and $cand.Range.ParagraphFormat.LeftIndent -eq 0) {
            $resourcesParaIdx = $k
            break
        }
    }
}

$p = $d.Paragraphs.Item($resourcesParaIdx)
$r = $p.Range
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:jc w:val="center"/><w:rPr><w:b/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="40"/><w:szCs w:val="40"/><w:u w:val="single"/></w:rPr><w:t>Resources</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Search for python for beginners by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Morsh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>youtube</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
[void]$r.InsertXML($xml)

$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$last.Range.ParagraphFormat.FirstLineIndent = 0

Write-Output "edit complete"
